# Auto-generated edit script to update cryptos list values
# Commit: Updated cryptos list on Thu Jun  8 12:26:38 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices) that Excel would otherwise
# auto-convert to real numbers; force the whole data range to Text format
# first so the assigned values remain literal strings exactly as scraped.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.418.67'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '1.843.71'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '264.48'
$ws.Range("E5").Value = '  -3.11%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.5184'
$ws.Range("E7").Value = '  -2.35%  '
$ws.Range("D8").Value = '0.3274'
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").Value = '0.06798'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = '18.78'
$ws.Range("E10").Value = '  -5.16%  '
$ws.Range("D11").Value = '0.7774'
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").Value = '0.07756'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").Value = '1.855.57'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '87.89'
$ws.Range("E14").Value = '  -2.45%  '
$ws.Range("D15").Value = '5.012'
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '13.90'
$ws.Range("E17").Value = '  -3.54%  '
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '0.000007944'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").Value = '26.428.95'
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").Value = '2.068.13'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").Value = '4.631'
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("E23").Value = '  -3.44%  '
$ws.Range("D24").Value = '5.988'
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").Value = '144.40'
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("D26").Value = '2.178'
$ws.Range("E26").Value = '  -8.44%  '
$ws.Range("D27").Value = '1.654'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '16.98'
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").Value = '111.96'
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").Value = '4.163'
$ws.Range("E30").Value = '  -3.47%  '
$ws.Range("D31").Value = '4.124'
$ws.Range("E31").Value = '  -4.06%  '
$ws.Range("D32").Value = '0.08711'
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("D33").Value = '0.04828'
$ws.Range("E33").Value = '  -2.00%  '
$ws.Range("D34").Value = '0.7245'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = '2.845'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").Value = '3.090'
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01780'
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.219'
$ws.Range("E39").Value = '  -4.07%  '
$ws.Range("D40").Value = '0.4859'
$ws.Range("E40").Value = '  -4.23%  '
$ws.Range("D41").Value = '0.9161'
$ws.Range("E41").Value = '  -2.12%  '
$ws.Range("D42").Value = '111.18'
$ws.Range("E42").Value = '  -4.17%  '
$ws.Range("D43").Value = '6.075'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '7.744'
$ws.Range("E45").Value = '  -3.29%  '
$ws.Range("D46").Value = '0.4170'
$ws.Range("E46").Value = '  -5.30%  '
$ws.Range("D47").Value = '0.05931'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").Value = '9.049'
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("D49").Value = '0.1246'
$ws.Range("E49").Value = '  -5.67%  '
$ws.Range("D50").Value = '34.98'
$ws.Range("E50").Value = '  -2.99%  '
$ws.Range("D51").Value = '0.8855'
$ws.Range("E51").Value = '  +0.82%  '

# Restore the original (default) cell style on column D now that the
# text values are set, so formatting matches the untouched workbook.
$ws.Range("D2:D51").Style = "Normal"
